# Airport Challenge - Domain Model.xlsx
# Tests added for returning plane object when land/takeoff
#
# - landPlane / takeOffPlane message cells now mention the "return:" path
#   alongside Log/@Error.
# - New "@Plane" return-type cells added in column L (rows 18 & 19) next to
#   those two message rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$quote = [char]0x201C
$oldFragment = "@Error " + $quote
$newFragment = "@Error / return: " + $quote

# - landPlane(@Plane) message cell
$landMsg = $ws.Range("D18").Value()
$ws.Range("D18").Value = $landMsg.Replace($oldFragment, $newFragment)

# - takeOffPlane(@Plane) message cell
$takeOffMsg = $ws.Range("D19").Value()
$ws.Range("D19").Value = $takeOffMsg.Replace($oldFragment, $newFragment)

# New "return value" column: both methods now return @Plane
$planeType = $ws.Range("D6").Value()
$ws.Range("L18").Value = $planeType
$ws.Range("L19").Value = $planeType

# Give the new column a sensible width (closest attainable value to 14.07 chars)
$ws.Columns.Item(11).ColumnWidth = 13.17

# Move the selection, matching where the author ended up editing
$ws.Range("J24").Select() | Out-Null
